# Actualización automática: inserts a new client row
# "JIMENEZ CORDERO WILLIAM GUSTAVO" (under asesor "OFICINA-CATAECSA")
# right before "LOZANO MOLINA TITO JERSON" in both the "VENTAS POR GRUPO"
# and "VENTA MENSUAL" sheets, pushing the remaining rows (and the trailing
# totals row) down by one.

$wb = $excel.ActiveWorkbook

$newAsesor  = "OFICINA-CATAECSA"
$newCliente = "JIMENEZ CORDERO WILLIAM GUSTAVO"
$insertAt   = 249

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A..R, data through row 288,
# totals-style row afterwards)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item($insertAt).Insert()

$ws1.Cells.Item($insertAt, 1).Value = $newAsesor
$ws1.Cells.Item($insertAt, 2).Value = $newCliente
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item($insertAt, $c).Value = 0
}

# The final summary row ("N de 287" counts) shifted from row 289 to 290;
# bump the "de 287" denominator to "de 288" to reflect the extra row.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(290, $c)
    $cell.Value = $cell.Value2 -replace "de 287", "de 288"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A..G, data through row 288, numeric
# totals row afterwards)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item($insertAt).Insert()

$ws2.Cells.Item($insertAt, 1).Value = $newAsesor
$ws2.Cells.Item($insertAt, 2).Value = $newCliente
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item($insertAt, $c).Value = 0
}
